$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text cell with new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.15 = 11791.91 pesos`n✅ 11791.91 pesos = 3.13 = 973.18 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the "tasas" sheet rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 317.6
$ws2.Range("O10").Value = 3745.11
$ws2.Range("N12").Value = 3762.28
$ws2.Range("O12").Value = 310.5
